# Refresh cryptos list (Price / Volume(1h) columns) to match the latest
# scrape. Row 14/15 also swap rank (WrappedliquidstakedEther2.0 <-> ShibaInu).
# Numeric-looking text in column D is prefixed with a leading apostrophe so
# Excel keeps it as literal text (e.g. "1.00") instead of coercing it to a
# number and dropping the trailing zero / formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.884.60"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "2.568.03"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'600.43"
$ws.Range("E5").Value = "  +1.89%  "
$ws.Range("D6").Value = "'178.49"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "2.568.31"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +11.78%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'0.345"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "'0.0000183"
$ws.Range("E14").Value = "  +5.88%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.006.44"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "'26.31"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "69.742.32"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("D18").Value = "2.575.92"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("D19").Value = "'7.78"
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'365.89"
$ws.Range("E21").Value = "  +3.66%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'70.82"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("E26").Value = "  +1.53%  "
$ws.Range("D27").Value = "'9.23"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "0.0₃0920"
$ws.Range("E30").Value = "  +0.21%  "
$ws.Range("D31").Value = "'514.61"
$ws.Range("E31").Value = "  +0.98%  "
$ws.Range("D32").Value = "'7.81"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'164.05"
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("E37").Value = "  -1.63%  "
$ws.Range("D38").Value = "'19.00"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("E39").Value = "  +1.40%  "
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("E44").Value = "  -1.39%  "
$ws.Range("E45").Value = "  -1.14%  "
$ws.Range("D46").Value = "'39.02"
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").Value = "'151.98"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("E48").Value = "  +1.60%  "
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("E51").Value = "  +1.64%  "
